# Update products.xlsx: Feuil1 sheet
# - Replace the "volume" column (C) values for rows 2-29 with new data
# - Move the active cell selection from B27 to D15

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

$volumeUpdates = @{
    2  = 50
    3  = 10
    4  = 10
    5  = 10
    6  = 30
    7  = 10
    8  = 10
    9  = 40
    10 = 100
    11 = 120
    12 = 10
    13 = 40
    14 = 20
    15 = 10
    16 = 10
    17 = 20
    18 = 10
    19 = 10
    20 = 2
    21 = 30
    22 = 40
    23 = 10
    24 = 10
    25 = 20
    26 = 1
    27 = 20
    28 = 20
    29 = 10
}

foreach ($row in $volumeUpdates.Keys) {
    $ws.Range("C$row").Value = $volumeUpdates[$row]
}

# Update the selected cell shown in the sheet view
$ws.Range("D15").Select()
